$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for new columns I and J ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style of the existing header cell (H1) onto the new header cells
# so they keep the same bold / bordered / centered look as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data values for columns I (I0) and J (IF), rows 2-84 ---
$iVals = @(2,7,10,3,6,1,3,9,9,9,5,7,6,7,8,7,8,7,6,7,7,8,8,8,8,7,6,7,8,2,8,8,6,6,1,7,6,8,4,6,7,9,7,7,6,8,8,6,8,4,7,8,6,7,7,9,7,7,7,8,7,7,5,6,6,4,9,5,7,5,6,6,6,8,5,5,5,4,5,5,6,8,4)
$jVals = @(3,7,11,4,7,1,4,9,9,9,6,7,6,7,8,7,8,7,7,8,7,8,8,8,8,8,7,7,8,4,8,8,7,6,2,7,6,8,5,7,7,9,8,7,6,8,8,7,8,4,7,8,7,7,7,9,7,7,8,8,7,7,5,6,6,5,9,6,7,6,7,7,7,8,6,7,5,4,5,6,6,8,4)

for ($r = 2; $r -le 84; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
